$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string edit: update source name for Instituto Nacional de Migración ---
$ws.Range("A9").Value = "Instituto Nacional de Migración (INM)"

# --- Row 10: new record "Comision Nacional de Bancos y Seguro (CNBS)" ---
$ws.Range("A10").Value = "Comision Nacional de Bancos y Seguro (CNBS)"
$ws.Range("C10").Value = "Trabajo"
$ws.Range("D10").Value = "Institución que por mandato constitucional tiene la responsabilidad de velar por la estabilidad y solvencia del sistema financiero y demás supervisados, su regulación, supervisión y control. Asimismo, vigilamos la transparencia y que se respeten los derechos de los usuarios financieros, así como coadyuvamos con el sistema de prevención y detección del lavado activos y financiamiento al terrorismo, y contribuimos a promover la educación e inclusión financiera, a fin de salvaguardar el interés público."
$ws.Range("E10").Value = "https://covid19honduras.org/?q=cnbs-22-3"
$ws.Range("F10").Value = "Las instituciones por la CNBS que realizan operaciones de crédito, podran otorgar periodos de gracia a los deudores que sean afectados por la reduccion de sus flujos de efectivo los cuales se podran otorgar hasta el 30 de junio de 2020."
$ws.Range("G10").Value = "https://covid19honduras.org/?q=cnbs-22-3"
$ws.Range("H10").Value = "21/3/2020"
$ws.Range("I10").Value = "22/3/2020"
$ws.Range("J10").Value = "Honduras"

$ws.Hyperlinks.Add($ws.Range("G10"), "https://covid19honduras.org/?q=cnbs-22-3")
$ws.Hyperlinks.Add($ws.Range("E10"), "https://covid19honduras.org/?q=cnbs-22-3")

# Re-apply the formatting of the row above (values/hyperlinks can reset a cell's
# style, so copy the reference row's formats last to land on the same style ids)
$ws.Range("A9:K9").Copy()
$ws.Range("A10:K10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows("10").RowHeight = 105

# --- Row 11: new record "Secretaria de Trabajo y Seguridad Social" ---
$ws.Range("A11").Value = "Secretaria de Trabajo y Seguridad Social"
$ws.Range("C11").Value = "Trabajo"
$ws.Range("G11").Value = "https://covid19honduras.org/?q=secretaria-de-trabajo"
$ws.Range("I11").Value = "26/3/2020"
$ws.Range("J11").Value = "Honduras"

$ws.Hyperlinks.Add($ws.Range("G11"), "https://covid19honduras.org/?q=secretaria-de-trabajo")

# Only the G column needs the hyperlink-style formatting here (other cells in
# row 11 already carry their target styles); copy from the row above.
$ws.Range("G9").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows("11").RowHeight = 45

# --- View state: selection moved as part of the edit session ---
$ws.Range("G11").Select()

Write-Output "Done"
